$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.448858737945557
$ws.Range("B1").Value = 2.38878870010376
$ws.Range("C1").Value = 2.807727098464966
$ws.Range("D1").Value = 3.22675347328186
$ws.Range("E1").Value = 1.858639478683472
